# Logged Week 15 and simulated Week 16
# Appends new per-play/per-game logged numbers to the running shared-string
# logs on the YDS and ST sheets, and updates the aggregate totals on the
# OFF, DEF, ST, TURNS and PEN sheets accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# YDS sheet: append newly logged play-by-play yardage values
# ---------------------------------------------------------------------
$ydsWs = $wb.Worksheets.Item("YDS")

$ydsWs.Range("B2").Value = $ydsWs.Range("B2").Value() + " 1 11 4 21 3 -2 0 7 8 -4 5 2 2 0 2 30 3 4 1 13 -1 1 7 1 7 8 5 7 5"
$ydsWs.Range("C2").Value = $ydsWs.Range("C2").Value() + " 8 2 16 6 3 0 2 2 26 3 4 1 -2 12 3 1 2 4 1 8 7 4 5 2 1"
$ydsWs.Range("B3").Value = $ydsWs.Range("B3").Value() + " 1 9 8 8 9 4 11 13 9 16 4 23 5 6 8 5 8"
$ydsWs.Range("C3").Value = $ydsWs.Range("C3").Value() + " 12 9 15 9 10 12 11 4 10 9 7 11 7 7 7 8 28 20"

# ---------------------------------------------------------------------
# OFF sheet: updated Home (row 2) / Road (row 3) season totals
# ---------------------------------------------------------------------
$offWs = $wb.Worksheets.Item("OFF")

$offWs.Range("C2").Value = 391
$offWs.Range("D2").Value = 15
$offWs.Range("E2").Value = 13
$offWs.Range("F2").Value = 89
$offWs.Range("G2").Value = 113
$offWs.Range("H2").Value = 10
$offWs.Range("I2").Value = 18
$offWs.Range("J2").Value = 63
$offWs.Range("N2").Value = 41
$offWs.Range("O2").Value = 55
$offWs.Range("P2").Value = 27

$offWs.Range("C3").Value = 312
$offWs.Range("E3").Value = 69
$offWs.Range("F3").Value = 215
$offWs.Range("G3").Value = 55
$offWs.Range("H3").Value = 47
$offWs.Range("I3").Value = 118
$offWs.Range("J3").Value = 83
$offWs.Range("L3").Value = 515
$offWs.Range("M3").Value = 331
$offWs.Range("Q3").Value = 953

# ---------------------------------------------------------------------
# DEF sheet: updated Home (row 2) / Road (row 3) season totals
# ---------------------------------------------------------------------
$defWs = $wb.Worksheets.Item("DEF")

$defWs.Range("C2").Value = 373
$defWs.Range("F2").Value = 104
$defWs.Range("G2").Value = 100
$defWs.Range("I2").Value = 11
$defWs.Range("J2").Value = 56
$defWs.Range("N2").Value = 36
$defWs.Range("O2").Value = 30
$defWs.Range("P2").Value = 17

$defWs.Range("B3").Value = 13
$defWs.Range("C3").Value = 299
$defWs.Range("D3").Value = 12
$defWs.Range("E3").Value = 65
$defWs.Range("F3").Value = 214
$defWs.Range("H3").Value = 56
$defWs.Range("I3").Value = 98
$defWs.Range("J3").Value = 103
$defWs.Range("L3").Value = 523
$defWs.Range("M3").Value = 356
$defWs.Range("Q3").Value = 914

# ---------------------------------------------------------------------
# ST sheet: updated totals (row 2) and appended per-return logs (rows 3-6)
# ---------------------------------------------------------------------
$stWs = $wb.Worksheets.Item("ST")

$stWs.Range("B2").Value = 142
$stWs.Range("D2").Value = 101
$stWs.Range("H2").Value = 6
$stWs.Range("I2").Value = 4

$stWs.Range("D3").Value = $stWs.Range("D3").Value() + " 47 54 41 40"
$stWs.Range("B4").Value = $stWs.Range("B4").Value() + " 49 47 50"
$stWs.Range("D4").Value = $stWs.Range("D4").Value() + " 17 14 0 0"
$stWs.Range("B5").Value = $stWs.Range("B5").Value() + " 13 8 28"
$stWs.Range("D5").Value = $stWs.Range("D5").Value() + " 0 0 0 0 16"
$stWs.Range("B6").Value = $stWs.Range("B6").Value() + " 21 10 18 21"

# ---------------------------------------------------------------------
# TURNS sheet: updated Road (row 3) totals
# ---------------------------------------------------------------------
$turnsWs = $wb.Worksheets.Item("TURNS")

$turnsWs.Range("B3").Value = 12
$turnsWs.Range("C3").Value = 9
$turnsWs.Range("D3").Value = 14
$turnsWs.Range("E3").Value = 18

# ---------------------------------------------------------------------
# PEN sheet: updated total (row 2)
# ---------------------------------------------------------------------
$penWs = $wb.Worksheets.Item("PEN")

$penWs.Range("D2").Value = 18
